# Applies the "Added most recent Amazon Order for heatshrink tubing to
# expense report" edit: row 12 (Programming Cable USB line item) is updated
# to reflect a newer Amazon order — an AmazonBasics USB cable replacing the
# previous StarTech cable entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update row 12 cell contents -------------------------------------------------

# Manufacturer
$ws.Range("C12").Value = "AmazonBasics"

# Vendor's Part #
$ws.Range("F12").Value = "B0718XZKWH"

# Manufacturer's Part #
$ws.Range("D12").Value = "3S44_24"

# Unit Price - now computed from the new order (24-pack at $47.99)
$ws.Range("G12").Formula = "=47.99/24"

# Description
$ws.Range("I12").Value = "USB 2.0 Cable - A-Male to Mini-B, 3 Feet (0.9 Meters)"

# Link - drop the old hyperlink and add the new order's link, mirroring how
# Excel recreates the hyperlink (new relationship id) when the link target
# changes.
$ws.Range("J12").Hyperlinks.Delete()
$ws.Range("J12").Value = "https://www.amazon.com/gp/product/B0718XZKWH/ref=oh_aui_detailpage_o03_s00?ie=UTF8&psc=1"
$ws.Hyperlinks.Add($ws.Range("J12"), "https://www.amazon.com/gp/product/B0718XZKWH/ref=oh_aui_detailpage_o03_s00?ie=UTF8&psc=1", "", "", "https://www.amazon.com/gp/product/B0718XZKWH/ref=oh_aui_detailpage_o03_s00?ie=UTF8&psc=1") | Out-Null
$ws.Range("J12").Style = "Hyperlink"

# --- Update view state: scrolled right with the new line item selected -----------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("J12").Select()

$wb.Save()
